# Update cryptos list data rows per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.131.03"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "3.882.24"
$ws.Range("E3").Value = "  -1.31%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "481.47"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").Value = "144.75"
$ws.Range("E6").Value = "  -2.92%  "

$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").Value = "0.737"
$ws.Range("E9").Value = "  +1.83%  "

$ws.Range("D10").Value = "0.183"
$ws.Range("E10").Value = "  +9.03%  "

$ws.Range("D11").Value = "0.0000357"
$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("D12").Value = "43.07"
$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("D13").Value = "10.47"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").Value = "4.497.67"
$ws.Range("E14").Value = "  -1.63%  "

$ws.Range("D15").Value = "3.864.97"
$ws.Range("E15").Value = "  -2.36%  "

$ws.Range("D16").Value = "14.18"
$ws.Range("E16").Value = "  -3.44%  "

$ws.Range("E17").Value = "  -0.64%  "

$ws.Range("D18").Value = "19.89"
$ws.Range("E18").Value = "  +0.46%  "

$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").Value = "68.159.58"
$ws.Range("E20").Value = "  -0.39%  "

$ws.Range("D21").Value = "'428.60"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("D22").Value = "3.52"
$ws.Range("E22").Value = "  +3.86%  "

$ws.Range("D23").Value = "14.75"
$ws.Range("E23").Value = "  +1.92%  "

$ws.Range("D24").Value = "89.06"
$ws.Range("E24").Value = "  +2.13%  "

$ws.Range("D25").Value = "12.01"
$ws.Range("E25").Value = "  +14.18%  "

$ws.Range("D26").Value = "3.65"
$ws.Range("E26").Value = "  +2.89%  "

$ws.Range("D27").Value = "'11.00"
$ws.Range("E27").Value = "  -1.49%  "

$ws.Range("D28").Value = "37.32"
$ws.Range("E28").Value = "  -2.42%  "

$ws.Range("E29").Value = "  -3.41%  "

$ws.Range("D30").Value = "709.97"
$ws.Range("E30").Value = "  -0.71%  "

$ws.Range("D31").Value = "13.45"
$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("D33").Value = "2.89"
$ws.Range("E33").Value = "  +2.04%  "

$ws.Range("E34").Value = "  +9.51%  "

$ws.Range("E35").Value = "  -0.97%  "

$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "60.96"
$ws.Range("E36").Value = "  +3.44%  "

$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "40.65"
$ws.Range("E37").Value = "  -2.64%  "

$ws.Range("D38").Value = "0.0502"
$ws.Range("E38").Value = "  +6.75%  "

$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.145"
$ws.Range("E40").Value = "  -4.24%  "

$ws.Range("D41").Value = "0.392"
$ws.Range("E41").Value = "  +14.43%  "

$ws.Range("D42").Value = "2.96"
$ws.Range("E42").Value = "  +3.91%  "

$ws.Range("D43").Value = "3.07"
$ws.Range("E43").Value = "  +3.10%  "

$ws.Range("D44").Value = "2.96"
$ws.Range("E44").Value = "  -2.38%  "

$ws.Range("E45").Value = "  +0.82%  "

$ws.Range("D46").Value = "3.36"
$ws.Range("E46").Value = "  +3.87%  "

$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("E48").Value = "  -1.45%  "

$ws.Range("E49").Value = "  -2.85%  "

$ws.Range("D50").Value = "143.61"
$ws.Range("E50").Value = "  -2.25%  "

$ws.Range("E51").Value = "  -2.19%  "
